$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.117504239082336
$ws.Range("B1").Value = 1.347274899482727
$ws.Range("C1").Value = 1.809825778007507
$ws.Range("D1").Value = 3.407470703125
$ws.Range("E1").Value = 1.99773120880127
